$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.8936481850813
$ws.Range("C2").Value = 10.59735749239121
$ws.Range("D2").Value = 8.166759778793239
$ws.Range("F2").Value = 44.00087492219443
$ws.Range("G2").Value = 3.713931810602639
$ws.Range("I2").Value = 28.88440754206507
$ws.Range("J2").Value = 11.52603766646551
$ws.Range("K2").Value = 14.88178688550446
$ws.Range("N2").Value = 21.68586532855986
$ws.Range("B3").Value = 13.68454498280087
$ws.Range("C3").Value = 10.43867091367731
$ws.Range("D3").Value = 8.109599212614709
$ws.Range("F3").Value = 43.93603146977713
$ws.Range("G3").Value = 3.716857207070902
$ws.Range("I3").Value = 28.92071921284065
$ws.Range("J3").Value = 11.51293872351883
$ws.Range("K3").Value = 14.74279332933433
$ws.Range("N3").Value = 21.74407199604368
$ws.Range("B4").Value = 13.55840901770215
$ws.Range("C4").Value = 10.34322086755793
$ws.Range("D4").Value = 8.076200178254076
$ws.Range("F4").Value = 43.9063962854762
$ws.Range("G4").Value = 3.718747075407678
$ws.Range("I4").Value = 28.94839315325247
$ws.Range("J4").Value = 11.50727054991812
$ws.Range("K4").Value = 14.66064865370474
$ws.Range("N4").Value = 21.78174008309983
$ws.Range("B5").Value = 13.50764310668389
$ws.Range("C5").Value = 10.30487423839076
$ws.Range("D5").Value = 8.063029121009951
$ws.Range("F5").Value = 43.89688380662415
$ws.Range("G5").Value = 3.719540847657932
$ws.Range("I5").Value = 28.96102003402915
$ws.Range("J5").Value = 11.505559173429
$ws.Range("K5").Value = 14.62801299815201
$ws.Range("N5").Value = 21.79757581992575
$ws.Range("B6").Value = 13.49925383689722
$ws.Range("C6").Value = 10.29854150386067
$ws.Range("D6").Value = 8.06086895826542
$ws.Range("F6").Value = 43.8954592223866
$ws.Range("G6").Value = 3.719674082909986
$ws.Range("I6").Value = 28.96319813410103
$ws.Range("J6").Value = 11.50531117636661
$ws.Range("K6").Value = 14.62264552788335
$ws.Range("N6").Value = 21.80023468506425
$ws.Range("B7").Value = 13.55772170770595
$ws.Range("C7").Value = 10.34270141767414
$ws.Range("D7").Value = 8.076020754159801
$ws.Range("F7").Value = 43.90625761066677
$ws.Range("G7").Value = 3.718757684689031
$ws.Range("I7").Value = 28.94855798380075
$ws.Range("J7").Value = 11.50724504514578
$ws.Range("K7").Value = 14.66020507619291
$ws.Range("N7").Value = 21.78195168205934
$ws.Range("B8").Value = 13.82112298801668
$ws.Range("C8").Value = 10.54226149486832
$ws.Range("D8").Value = 8.146706105228017
$ws.Range("F8").Value = 43.97640758638906
$ws.Range("G8").Value = 3.714921097718156
$ws.Range("I8").Value = 28.89580984542014
$ws.Range("J8").Value = 11.52102921359962
$ws.Range("K8").Value = 14.83321957424218
$ws.Range("N8").Value = 21.70553501345188
$ws.Range("B9").Value = 14.35249076596163
$ws.Range("C9").Value = 10.94705882343282
$ws.Range("D9").Value = 8.298212019862747
$ws.Range("F9").Value = 44.1944502990545
$ws.Range("G9").Value = 3.70813693947324
$ws.Range("I9").Value = 28.83516633386812
$ws.Range("J9").Value = 11.56682667001515
$ws.Range("K9").Value = 15.19626834355792
$ws.Range("N9").Value = 21.57095812667518
$ws.Range("B10").Value = 14.74785463664231
$ws.Range("C10").Value = 11.24959701496158
$ws.Range("D10").Value = 8.416563668046926
$ws.Range("F10").Value = 44.40319614351255
$ws.Range("G10").Value = 3.703598046345637
$ws.Range("I10").Value = 28.81684888471067
$ws.Range("J10").Value = 11.61179045826198
$ws.Range("K10").Value = 15.4752492467921
$ws.Range("N10").Value = 21.48135518972699
$ws.Range("B11").Value = 14.9279460075841
$ws.Range("C11").Value = 11.38770023459126
$ws.Range("D11").Value = 8.471748409061533
$ws.Range("F11").Value = 44.5085486139165
$ws.Range("G11").Value = 3.701628773244371
$ws.Range("I11").Value = 28.81423748687194
$ws.Range("J11").Value = 11.63466778198404
$ws.Range("K11").Value = 15.60432954337271
$ws.Range("N11").Value = 21.44259693250275
$ws.Range("B12").Value = 14.99610754469047
$ws.Range("C12").Value = 11.44001256828957
$ws.Range("D12").Value = 8.492822837615549
$ws.Range("F12").Value = 44.5499207556925
$ws.Range("G12").Value = 3.700896705862854
$ws.Range("I12").Value = 28.81407260364956
$ws.Range("J12").Value = 11.6436755690907
$ws.Range("K12").Value = 15.65347880149712
$ws.Range("N12").Value = 21.42820760007767
$ws.Range("B13").Value = 14.98143038339692
$ws.Range("C13").Value = 11.42874630234517
$ws.Range("D13").Value = 8.488276466539061
$ws.Range("F13").Value = 44.5409451083468
$ws.Range("G13").Value = 3.701053763662073
$ws.Range("I13").Value = 28.81407145066557
$ws.Range("J13").Value = 11.64172031774008
$ws.Range("K13").Value = 15.64288233150002
$ws.Range("N13").Value = 21.43129381605376
$ws.Range("B14").Value = 14.93355479263263
$ws.Range("C14").Value = 11.39200398616394
$ws.Range("D14").Value = 8.473478759114579
$ws.Range("F14").Value = 44.51192281971076
$ws.Range("G14").Value = 3.701568272427246
$ws.Range("I14").Value = 28.8142074018325
$ws.Range("J14").Value = 11.63540197412038
$ws.Range("K14").Value = 15.60836792943189
$ws.Range("N14").Value = 21.44140735210341
$ws.Range("B15").Value = 14.9042230096029
$ws.Range("C15").Value = 11.3694987478572
$ws.Range("D15").Value = 8.464437317317653
$ws.Range("F15").Value = 44.49433769572548
$ws.Range("G15").Value = 3.701885199695149
$ws.Range("I15").Value = 28.81439801555057
$ws.Range("J15").Value = 11.63157656904352
$ws.Range("K15").Value = 15.58726066841749
$ws.Range("N15").Value = 21.44763962504412
$ws.Range("B16").Value = 14.73608405381356
$ws.Range("C16").Value = 11.24057668359461
$ws.Range("D16").Value = 8.412982905672866
$ws.Range("F16").Value = 44.39651883519506
$ws.Range("G16").Value = 3.703728657836935
$ws.Range("I16").Value = 28.81713479161764
$ws.Range("J16").Value = 11.61034381490316
$ws.Range("K16").Value = 15.46685356905479
$ws.Range("N16").Value = 21.4839283960042
$ws.Range("B17").Value = 14.63294722738613
$ws.Range("C17").Value = 11.16157122236984
$ws.Range("D17").Value = 8.381750462154086
$ws.Range("F17").Value = 44.33916055494909
$ws.Range("G17").Value = 3.70488396184833
$ws.Range("I17").Value = 28.82028011728188
$ws.Range("J17").Value = 11.5979362472357
$ws.Range("K17").Value = 15.39351325233438
$ws.Range("N17").Value = 21.50670300318978
$ws.Range("B18").Value = 14.57365211878851
$ws.Range("C18").Value = 11.11617736129914
$ws.Range("D18").Value = 8.363914078346813
$ws.Range("F18").Value = 44.30714897663419
$ws.Range("G18").Value = 3.705557454965867
$ws.Range("I18").Value = 28.82262769006563
$ws.Range("J18").Value = 11.59102806702357
$ws.Range("K18").Value = 15.35153671445033
$ws.Range("N18").Value = 21.51999086757519
$ws.Range("B19").Value = 14.55358251010921
$ws.Range("C19").Value = 11.10081766011335
$ws.Range("D19").Value = 8.357897411051422
$ws.Range("F19").Value = 44.29647910203041
$ws.Range("G19").Value = 3.705787035074356
$ws.Range("I19").Value = 28.82351497414608
$ws.Range("J19").Value = 11.58872839965987
$ws.Range("K19").Value = 15.33736101113617
$ws.Range("N19").Value = 21.52452229967403
$ws.Range("B20").Value = 14.64392408933576
$ws.Range("C20").Value = 11.16997691013793
$ws.Range("D20").Value = 8.385062111311541
$ws.Range("F20").Value = 44.34516519557394
$ws.Range("G20").Value = 3.704760047608476
$ws.Range("I20").Value = 28.81988955227999
$ws.Range("J20").Value = 11.59923344943152
$ws.Range("K20").Value = 15.40129935308258
$ws.Range("N20").Value = 21.50425909986247
$ws.Range("B21").Value = 14.947618510589
$ws.Range("C21").Value = 11.40279608499694
$ws.Range("D21").Value = 8.477820528159365
$ws.Range("F21").Value = 44.52040741804939
$ws.Range("G21").Value = 3.70141677879578
$ws.Range("I21").Value = 28.81414509846302
$ws.Range("J21").Value = 11.63724850361261
$ws.Range("K21").Value = 15.61849867702082
$ws.Range("N21").Value = 21.43842895958375
$ws.Range("B22").Value = 15.14586489209724
$ws.Range("C22").Value = 11.55502462964447
$ws.Range("D22").Value = 8.539468403488398
$ws.Range("F22").Value = 44.64354067342819
$ws.Range("G22").Value = 3.699311307107043
$ws.Range("I22").Value = 28.8151939528233
$ws.Range("J22").Value = 11.66410032358261
$ws.Range("K22").Value = 15.76200114138223
$ws.Range("N22").Value = 21.39708122201612
$ws.Range("B23").Value = 15.04010095955801
$ws.Range("C23").Value = 11.47378827561302
$ws.Range("D23").Value = 8.506477470922897
$ws.Range("F23").Value = 44.57704119001085
$ws.Range("G23").Value = 3.700427784020862
$ws.Range("I23").Value = 28.81419435029336
$ws.Range("J23").Value = 11.6495867400066
$ws.Range("K23").Value = 15.68528347144327
$ws.Range("N23").Value = 21.41899604964863
$ws.Range("B24").Value = 14.63896144584604
$ws.Range("C24").Value = 11.16617660918601
$ws.Range("D24").Value = 8.383564541065299
$ws.Range("F24").Value = 44.34244749185566
$ws.Range("G24").Value = 3.704816040280354
$ws.Range("I24").Value = 28.8200644468447
$ws.Range("J24").Value = 11.59864628289084
$ws.Range("K24").Value = 15.39777866869493
$ws.Range("N24").Value = 21.5053633827709
$ws.Range("B25").Value = 14.20757656173074
$ws.Range("C25").Value = 10.83642745644794
$ws.Range("D25").Value = 8.255927800820615
$ws.Range("F25").Value = 44.12689240913951
$ws.Range("G25").Value = 3.709893629872627
$ws.Range("I25").Value = 28.84697372500909
$ws.Range("J25").Value = 11.55243928125746
$ws.Range("K25").Value = 15.09573681659627
$ws.Range("N25").Value = 21.60573376568861
